$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before row 75, pushing the "End with ..." summary block
# (old rows 75-78) down to rows 80-83, and leaving a fresh blank block at
# rows 75-79 for the new "Middle with ..." summary rows.
$ws.Range("A75:A79").EntireRow.Insert()

# --- New "Middle with ..." row labels (written first to match shared-string order) ---
$ws.Range("B75").Value = "Middle with A"
$ws.Range("B76").Value = "Middle with T"
$ws.Range("B77").Value = "Middle with G"
$ws.Range("B78").Value = "Middle with C"

# --- New header labels for the Pgene+Ngene / ProteinPGene+ProteinNGene columns ---
$ws.Range("I69").Value = "Pgene+Ngene"
$ws.Range("J69").Value = "ProteinPGene+ProteinNGene"

# --- New nucleotide labels for the grand-total block ---
$ws.Range("H85").Value = "A"
$ws.Range("H87").Value = "G"
$ws.Range("H86").Value = "T"
$ws.Range("H88").Value = "C"

# --- New column J formulas for the existing "Start with ..." rows (70-73) ---
$ws.Range("J70").Formula = "=SUM(C70+F70)"
$ws.Range("J71").Formula = "=SUM(C71+F71)"
$ws.Range("J72").Formula = "=SUM(C72+F72)"
$ws.Range("J73").Formula = "=SUM(C73+F73)"

# --- New "Middle with ..." summary block formulas (rows 75-78) ---
$ws.Range("C75").Formula = "=SUM(C3,C4,C5,C6,C19,C20,C21,C22,C35,C36,C37,C38,C51,C52,C53,C54)"
$ws.Range("D75").Formula = "=SUM(D3,D4,D5,D6,D19,D20,D21,D22,D35,D36,D37,D38,D51,D52,D53,D54)"
$ws.Range("E75").Formula = "=SUM(E3,E4,E5,E6,E19,E20,E21,E22,E35,E36,E37,E38,E51,E52,E53,E54)"
$ws.Range("F75").Formula = "=SUM(F3,F4,F5,F6,F19,F20,F21,F22,F35,F36,F37,F38,F51,F52,F53,F54)"
$ws.Range("G75").Formula = "=SUM(G3,G4,G5,G6,G19,G20,G21,G22,G35,G36,G37,G38,G51,G52,G53,G54)"
$ws.Range("H75").Formula = "=SUM(H3,H4,H5,H6,H19,H20,H21,H22,H35,H36,H37,H38,H51,H52,H53,H54)"
$ws.Range("I75").Formula = "=SUM(D75,H75)"
$ws.Range("J75").Formula = "=SUM(C75+F75)"

$ws.Range("C76").Formula = "=SUM(C7,C8,C9,C10,C23,C24,C25,C26,C39,C40,C41,C42,C55,C56,C57,C58)"
$ws.Range("D76").Formula = "=SUM(D7,D8,D9,D10,D23,D24,D25,D26,D39,D40,D41,D42,D55,D56,D57,D58)"
$ws.Range("E76").Formula = "=SUM(E7,E8,E9,E10,E23,E24,E25,E26,E39,E40,E41,E42,E55,E56,E57,E58)"
$ws.Range("F76").Formula = "=SUM(F7,F8,F9,F10,F23,F24,F25,F26,F39,F40,F41,F42,F55,F56,F57,F58)"
$ws.Range("G76").Formula = "=SUM(G7,G8,G9,G10,G23,G24,G25,G26,G39,G40,G41,G42,G55,G56,G57,G58)"
$ws.Range("H76").Formula = "=SUM(H7,H8,H9,H10,H23,H24,H25,H26,H39,H40,H41,H42,H55,H56,H57,H58)"
$ws.Range("I76").Formula = "=SUM(D76,H76)"
$ws.Range("J76").Formula = "=SUM(C76+F76)"

$ws.Range("C77").Formula = "=SUM(C11,C12,C13,C14,C27,C28,C29,C30,C43,C44,C45,C46,C59,C60,C61,C62)"
$ws.Range("D77").Formula = "=SUM(D11,D12,D13,D14,D27,D28,D29,D30,D43,D44,D45,D46,D59,D60,D61,D62)"
$ws.Range("E77").Formula = "=SUM(E11,E12,E13,E14,E27,E28,E29,E30,E43,E44,E45,E46,E59,E60,E61,E62)"
$ws.Range("F77").Formula = "=SUM(F11,F12,F13,F14,F27,F28,F29,F30,F43,F44,F45,F46,F59,F60,F61,F62)"
$ws.Range("G77").Formula = "=SUM(G11,G12,G13,G14,G27,G28,G29,G30,G43,G44,G45,G46,G59,G60,G61,G62)"
$ws.Range("H77").Formula = "=SUM(H11,H12,H13,H14,H27,H28,H29,H30,H43,H44,H45,H46,H59,H60,H61,H62)"
$ws.Range("I77").Formula = "=SUM(D77,H77)"
$ws.Range("J77").Formula = "=SUM(C77+F77)"

$ws.Range("C78").Formula = "=SUM(C15,C16,C17,C18,C31,C32,C33,C34,C47,C48,C49,C50,C63,C64,C65,C66)"
$ws.Range("D78").Formula = "=SUM(D15,D16,D17,D18,D31,D32,D33,D34,D47,D48,D49,D50,D63,D64,D65,D66)"
$ws.Range("E78").Formula = "=SUM(E15,E16,E17,E18,E31,E32,E33,E34,E47,E48,E49,E50,E63,E64,E65,E66)"
$ws.Range("F78").Formula = "=SUM(F15,F16,F17,F18,F31,F32,F33,F34,F47,F48,F49,F50,F63,F64,F65,F66)"
$ws.Range("G78").Formula = "=SUM(G15,G16,G17,G18,G31,G32,G33,G34,G47,G48,G49,G50,G63,G64,G65,G66)"
$ws.Range("H78").Formula = "=SUM(H15,H16,H17,H18,H31,H32,H33,H34,H47,H48,H49,H50,H63,H64,H65,H66)"
$ws.Range("I78").Formula = "=SUM(D78,H78)"
$ws.Range("J78").Formula = "=SUM(C78+F78)"

# --- New column J formulas for the "End with ..." rows (now at 80-83) ---
$ws.Range("J80").Formula = "=SUM(C80+F80)"
$ws.Range("J81").Formula = "=SUM(C81+F81)"
$ws.Range("J82").Formula = "=SUM(C82+F82)"
$ws.Range("J83").Formula = "=SUM(C83+F83)"

# --- New grand-total block (rows 85-88) ---
$ws.Range("I85").Formula = "=SUM(I70+I75+I80)"
$ws.Range("J85").Formula = "=SUM(J70+J75+J80)"

$ws.Range("I86").Formula = "=SUM(I71+I76+I81)"
$ws.Range("J86").Formula = "=SUM(J71+J76+J81)"

$ws.Range("I87").Formula = "=SUM(I72+I77+I82)"
$ws.Range("J87").Formula = "=SUM(J72+J77+J82)"

$ws.Range("I88").Formula = "=SUM(I73+I78+I83)"
$ws.Range("J88").Formula = "=SUM(J73+J78+J83)"

# --- View state: scrolled down to keep the bottom of the expanded table visible ---
$excel.ActiveWindow.ScrollRow = 67
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L88").Select()

# --- Recalculate everything ---
$excel.CalculateFull()
